$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.325.39'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '1.843.79'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9975'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.60'
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6266'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9984'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07493'
$ws.Range("E8").Value = '  -1.73%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2898'
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.41'
$ws.Range("E10").Value = '  -1.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07727'
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("D12").Value = '1.843.94'
$ws.Range("E12").Value = '  -2.34%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.983'
$ws.Range("E13").Value = '  -0.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6802'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001048'
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.01'
$ws.Range("E16").Value = '  -1.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.178'
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("D18").Value = '29.386.20'
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '229.18'
$ws.Range("E19").Value = '  +0.60%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.31'
$ws.Range("E20").Value = '  -0.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9985'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.475'
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9961'
$ws.Range("E23").Value = '  -0.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.56'
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.416'
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1370'
$ws.Range("E26").Value = '  -1.14%  '
$ws.Range("E27").Value = '  -0.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06381'
$ws.Range("E28").Value = '  +13.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.405'
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.477'
$ws.Range("E30").Value = '  +1.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.091'
$ws.Range("E31").Value = '  -0.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.093'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.831'
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("E34").Value = '  -1.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6965'
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("E36").Value = '  -0.37%  '
$ws.Range("D37").Value = '1.267.31'
$ws.Range("E37").Value = '  +3.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.841'
$ws.Range("E38").Value = '  +4.51%  '
$ws.Range("E39").Value = '  +1.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.720'
$ws.Range("E40").Value = '  +5.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9144'
$ws.Range("E41").Value = '  +1.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9977'
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("D43").Value = '2.005.80'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '66.11'
$ws.Range("E45").Value = '  +0.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.074'
$ws.Range("E46").Value = '  -1.75%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.723'
$ws.Range("E47").Value = '  +2.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1164'
$ws.Range("E48").Value = '  +1.88%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.994'
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.00000000114'
$ws.Range("E51").Value = '  -0.92%  '
